$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 17.48085021972656
$ws.Range("D2").Value = 123

$ws.Range("C3").Value = 16.82925224304199
$ws.Range("D3").Value = 123

$ws.Range("C4").Value = 17.46201515197754
$ws.Range("D4").Value = 174

$ws.Range("C5").Value = 17.36998558044434
$ws.Range("D5").Value = 123

$ws.Range("C6").Value = 18.92304420471191
$ws.Range("D6").Value = 123
